# Apply the latest crypto price/volume snapshot to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a cell reference plus its new text value. The "Price" column
# (D) stores values as plain text (they use "." as a thousands separator and
# would otherwise be misread as numbers), so cells whose new value would parse
# as a number are temporarily switched to a text format before assignment and
# restored to the default style afterwards, exactly like the rest of the sheet.
$changes = @(
    @{ Cell = 'D2'; Value = '28.310.57'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  -1.69%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.551.80'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  -1.49%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  -0.01%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '209.83'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  -1.64%  '; ForceText = $false }
    @{ Cell = 'E6'; Value = '  -1.86%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  -0.01%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '23.74'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  -2.07%  '; ForceText = $false }
    @{ Cell = 'E9'; Value = '  -2.10%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.0583'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  -1.68%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.0891'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  +0.25%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '1.774.48'; ForceText = $false }
    @{ Cell = 'E12'; Value = '  -1.42%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '1.559.31'; ForceText = $false }
    @{ Cell = 'E13'; Value = '  -0.83%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '28.284.24'; ForceText = $false }
    @{ Cell = 'E14'; Value = '  -1.76%  '; ForceText = $false }
    @{ Cell = 'E15'; Value = '  -1.95%  '; ForceText = $false }
    @{ Cell = 'E16'; Value = '  -2.53%  '; ForceText = $false }
    @{ Cell = 'E17'; Value = '  -3.15%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '227.88'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  -1.14%  '; ForceText = $false }
    @{ Cell = 'E19'; Value = '  -1.05%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '0.0₃0675'; ForceText = $false }
    @{ Cell = 'E20'; Value = '  -2.86%  '; ForceText = $false }
    @{ Cell = 'E22'; Value = '  +0.38%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '8.90'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  -3.99%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '151.46'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -0.34%  '; ForceText = $false }
    @{ Cell = 'E27'; Value = '  -1.27%  '; ForceText = $false }
    @{ Cell = 'E28'; Value = '  -0.06%  '; ForceText = $false }
    @{ Cell = 'E29'; Value = '  -3.33%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '0.0467'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  -3.71%  '; ForceText = $false }
    @{ Cell = 'E31'; Value = '  -4.68%  '; ForceText = $false }
    @{ Cell = 'E32'; Value = '  -1.58%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '1.387.91'; ForceText = $false }
    @{ Cell = 'E33'; Value = '  -0.16%  '; ForceText = $false }
    @{ Cell = 'E34'; Value = '  -2.91%  '; ForceText = $false }
    @{ Cell = 'E35'; Value = '  +1.06%  '; ForceText = $false }
    @{ Cell = 'E36'; Value = '  -3.68%  '; ForceText = $false }
    @{ Cell = 'E37'; Value = '  -1.19%  '; ForceText = $false }
    @{ Cell = 'E38'; Value = '  -1.09%  '; ForceText = $false }
    @{ Cell = 'E39'; Value = '  -3.23%  '; ForceText = $false }
    @{ Cell = 'E40'; Value = '  +1.16%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '0.510'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -3.22%  '; ForceText = $false }
    @{ Cell = 'E42'; Value = '  -0.07%  '; ForceText = $false }
    @{ Cell = 'E44'; Value = '  -0.92%  '; ForceText = $false }
    @{ Cell = 'E45'; Value = '  -2.53%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '61.92'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  -2.35%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '1.687.89'; ForceText = $false }
    @{ Cell = 'E47'; Value = '  -1.43%  '; ForceText = $false }
    @{ Cell = 'E48'; Value = '  -6.21%  '; ForceText = $false }
    @{ Cell = 'E49'; Value = '  -1.19%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '42.35'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  +5.84%  '; ForceText = $false }
    @{ Cell = 'E51'; Value = '  +0.38%  '; ForceText = $false }
)

foreach ($chg in $changes) {
    $cell = $ws.Range($chg.Cell)
    if ($chg.ForceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $chg.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $chg.Value
    }
}

Write-Output "Applied $($changes.Count) cell updates to Sheet1"
